$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI ligand-receptor statistics ("Natmi following Dr Hou advice")
# Ligand/Receptor-expressing cell counts (E, K) increased from 1 to 3 for rows 2-17,
# with all derived expression/specificity/edge-weight statistics recalculated accordingly.
# Columns F and L (detection rate) are unchanged.

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 18.04537966666667
$ws.Range("H2").Value = 54.13613900000001
$ws.Range("I2").Value = 0.6797959733292525
$ws.Range("J2").Value = 0.6797959733292525
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 15.35884066666667
$ws.Range("N2").Value = 46.076522
$ws.Range("O2").Value = 0.1012042817263867
$ws.Range("P2").Value = 0.1012042817263867
$ws.Range("Q2").Value = 277.1561110698398
$ws.Range("R2").Value = 2494.404999628558
$ws.Range("S2").Value = 0.0687982632012769
$ws.Range("T2").Value = 0.0687982632012769

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 18.04537966666667
$ws.Range("H3").Value = 54.13613900000001
$ws.Range("I3").Value = 0.6797959733292525
$ws.Range("J3").Value = 0.6797959733292525
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 50.59256466666667
$ws.Range("N3").Value = 151.777694
$ws.Range("O3").Value = 0.3333704853712116
$ws.Range("P3").Value = 0.3333704853712116
$ws.Range("Q3").Value = 912.9620377203853
$ws.Range("R3").Value = 8216.658339483467
$ws.Range("S3").Value = 0.2266239135821681
$ws.Range("T3").Value = 0.2266239135821681

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 18.04537966666667
$ws.Range("H4").Value = 54.13613900000001
$ws.Range("I4").Value = 0.6797959733292525
$ws.Range("J4").Value = 0.6797959733292525
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 60.37715666666667
$ws.Range("N4").Value = 181.13147
$ws.Range("O4").Value = 0.397844271305776
$ws.Range("P4").Value = 0.397844271305776
$ws.Range("Q4").Value = 1089.528715243815
$ws.Range("R4").Value = 9805.758437194332
$ws.Range("S4").Value = 0.2704529336457772
$ws.Range("T4").Value = 0.2704529336457772

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 18.04537966666667
$ws.Range("H5").Value = 54.13613900000001
$ws.Range("I5").Value = 0.6797959733292525
$ws.Range("J5").Value = 0.6797959733292525
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 25.43221733333333
$ws.Range("N5").Value = 76.29665199999999
$ws.Range("O5").Value = 0.1675809615966257
$ws.Range("P5").Value = 0.1675809615966258
$ws.Range("Q5").Value = 458.9340175451809
$ws.Range("R5").Value = 4130.406157906628
$ws.Range("S5").Value = 0.1139208629000303
$ws.Range("T5").Value = 0.1139208629000303

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.6001993333333334
$ws.Range("H6").Value = 1.800598
$ws.Range("I6").Value = 0.02261039099934159
$ws.Range("J6").Value = 0.02261039099934159
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 15.35884066666667
$ws.Range("N6").Value = 46.076522
$ws.Range("O6").Value = 0.1012042817263867
$ws.Range("P6").Value = 0.1012042817263867
$ws.Range("Q6").Value = 9.218365928906223
$ws.Range("R6").Value = 82.965293360156
$ws.Range("S6").Value = 0.002288268380641124
$ws.Range("T6").Value = 0.002288268380641124

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.6001993333333334
$ws.Range("H7").Value = 1.800598
$ws.Range("I7").Value = 0.02261039099934159
$ws.Range("J7").Value = 0.02261039099934159
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 50.59256466666667
$ws.Range("N7").Value = 151.777694
$ws.Range("O7").Value = 0.3333704853712116
$ws.Range("P7").Value = 0.3333704853712116
$ws.Range("Q7").Value = 30.36562358455689
$ws.Range("R7").Value = 273.290612261012
$ws.Range("S7").Value = 0.00753763702188338
$ws.Range("T7").Value = 0.00753763702188338

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.6001993333333334
$ws.Range("H8").Value = 1.800598
$ws.Range("I8").Value = 0.02261039099934159
$ws.Range("J8").Value = 0.02261039099934159
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 60.37715666666667
$ws.Range("N8").Value = 181.13147
$ws.Range("O8").Value = 0.397844271305776
$ws.Range("P8").Value = 0.397844271305776
$ws.Range("Q8").Value = 36.23832917989557
$ws.Range("R8").Value = 326.14496261906
$ws.Range("S8").Value = 0.008995414531071732
$ws.Range("T8").Value = 0.008995414531071732

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.6001993333333334
$ws.Range("H9").Value = 1.800598
$ws.Range("I9").Value = 0.02261039099934159
$ws.Range("J9").Value = 0.02261039099934159
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 25.43221733333333
$ws.Range("N9").Value = 76.29665199999999
$ws.Range("O9").Value = 0.1675809615966257
$ws.Range("P9").Value = 0.1675809615966258
$ws.Range("Q9").Value = 15.26439988865511
$ws.Range("R9").Value = 137.379598997896
$ws.Range("S9").Value = 0.003789071065745356
$ws.Range("T9").Value = 0.003789071065745357

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 4.690054333333333
$ws.Range("H10").Value = 14.070163
$ws.Range("I10").Value = 0.1766812397072912
$ws.Range("J10").Value = 0.1766812397072912
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 15.35884066666667
$ws.Range("N10").Value = 46.076522
$ws.Range("O10").Value = 0.1012042817263867
$ws.Range("P10").Value = 0.1012042817263867
$ws.Range("Q10").Value = 72.03379722367622
$ws.Range("R10").Value = 648.304175013086
$ws.Range("S10").Value = 0.01788089795910395
$ws.Range("T10").Value = 0.01788089795910395

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 4.690054333333333
$ws.Range("H11").Value = 14.070163
$ws.Range("I11").Value = 0.1766812397072912
$ws.Range("J11").Value = 0.1766812397072912
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 50.59256466666667
$ws.Range("N11").Value = 151.777694
$ws.Range("O11").Value = 0.3333704853712116
$ws.Range("P11").Value = 0.3333704853712116
$ws.Range("Q11").Value = 237.2818771493469
$ws.Range("R11").Value = 2135.536894344122
$ws.Range("S11").Value = 0.05890031063720703
$ws.Range("T11").Value = 0.05890031063720704

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 4.690054333333333
$ws.Range("H12").Value = 14.070163
$ws.Range("I12").Value = 0.1766812397072912
$ws.Range("J12").Value = 0.1766812397072912
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 60.37715666666667
$ws.Range("N12").Value = 181.13147
$ws.Range("O12").Value = 0.397844271305776
$ws.Range("P12").Value = 0.397844271305776
$ws.Range("Q12").Value = 283.1721452588456
$ws.Range("R12").Value = 2548.54930732961
$ws.Range("S12").Value = 0.07029161906474839
$ws.Range("T12").Value = 0.0702916190647484

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 4.690054333333333
$ws.Range("H13").Value = 14.070163
$ws.Range("I13").Value = 0.1766812397072912
$ws.Range("J13").Value = 0.1766812397072912
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 25.43221733333333
$ws.Range("N13").Value = 76.29665199999999
$ws.Range("O13").Value = 0.1675809615966257
$ws.Range("P13").Value = 0.1675809615966258
$ws.Range("Q13").Value = 119.2784811104751
$ws.Range("R13").Value = 1073.506329994276
$ws.Range("S13").Value = 0.02960841204623179
$ws.Range("T13").Value = 0.0296084120462318

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 3.209654333333333
$ws.Range("H14").Value = 9.628962999999999
$ws.Range("I14").Value = 0.1209123959641148
$ws.Range("J14").Value = 0.1209123959641148
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 15.35884066666667
$ws.Range("N14").Value = 46.076522
$ws.Range("O14").Value = 0.1012042817263867
$ws.Range("P14").Value = 0.1012042817263867
$ws.Range("Q14").Value = 49.29656950074288
$ws.Range("R14").Value = 443.6691255066859
$ws.Range("S14").Value = 0.01223685218536469
$ws.Range("T14").Value = 0.01223685218536469

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 3.209654333333333
$ws.Range("H15").Value = 9.628962999999999
$ws.Range("I15").Value = 0.1209123959641148
$ws.Range("J15").Value = 0.1209123959641148
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 50.59256466666667
$ws.Range("N15").Value = 151.777694
$ws.Range("O15").Value = 0.3333704853712116
$ws.Range("P15").Value = 0.3333704853712116
$ws.Range("Q15").Value = 162.3846444168135
$ws.Range("R15").Value = 1461.461799751322
$ws.Range("S15").Value = 0.04030862412995307
$ws.Range("T15").Value = 0.04030862412995307

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 3.209654333333333
$ws.Range("H16").Value = 9.628962999999999
$ws.Range("I16").Value = 0.1209123959641148
$ws.Range("J16").Value = 0.1209123959641148
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 60.37715666666667
$ws.Range("N16").Value = 181.13147
$ws.Range("O16").Value = 0.397844271305776
$ws.Range("P16").Value = 0.397844271305776
$ws.Range("Q16").Value = 193.7898025295122
$ws.Range("R16").Value = 1744.10822276561
$ws.Range("S16").Value = 0.0481043040641787
$ws.Range("T16").Value = 0.0481043040641787

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 3.209654333333333
$ws.Range("H17").Value = 9.628962999999999
$ws.Range("I17").Value = 0.1209123959641148
$ws.Range("J17").Value = 0.1209123959641148
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 25.43221733333333
$ws.Range("N17").Value = 76.29665199999999
$ws.Range("O17").Value = 0.1675809615966257
$ws.Range("P17").Value = 0.1675809615966258
$ws.Range("Q17").Value = 81.62862657020843
$ws.Range("R17").Value = 734.6576391318758
$ws.Range("S17").Value = 0.02026261558461833
$ws.Range("T17").Value = 0.02026261558461833

Write-Output "Updated rows 2-17 columns E,G,H,I,J,K,M,N,O,P,Q,R,S,T"
